# This workbook holds NATMI ligand-receptor signalling output (Ptn-Ptprz1).
# The underlying analysis script was re-run against updated TPM expression
# values, so every computed metric column (E:T) is refreshed in place below,
# cell by cell, while the categorical columns (A:D) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2.0
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2563003333333334
$ws.Range("H2").Value = 0.7689010000000001
$ws.Range("I2").Value = 0.02986826554325775
$ws.Range("J2").Value = 0.02986826554325775
$ws.Range("K2").Value = 1.0
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04936366666666667
$ws.Range("N2").Value = 0.148091
$ws.Range("O2").Value = 0.04616170608573571
$ws.Range("P2").Value = 0.0461617060857357
$ws.Range("Q2").Value = 0.01265192422122222
$ws.Range("R2").Value = 0.113867317991
$ws.Range("S2").Value = 0.001378770095298571
$ws.Range("T2").Value = 0.001378770095298571

# Row 3
$ws.Range("E3").Value = 2.0
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2563003333333334
$ws.Range("H3").Value = 0.7689010000000001
$ws.Range("I3").Value = 0.02986826554325775
$ws.Range("J3").Value = 0.02986826554325775
$ws.Range("O3").Value = 0.008057748967298944
$ws.Range("P3").Value = 0.008057748967298944
$ws.Range("Q3").Value = 0.002208454538888889
$ws.Range("R3").Value = 0.01987609085
$ws.Range("S3").Value = 0.0002406709858361958
$ws.Range("T3").Value = 0.0002406709858361957

# Row 4
$ws.Range("E4").Value = 2.0
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2563003333333334
$ws.Range("H4").Value = 0.7689010000000001
$ws.Range("I4").Value = 0.02986826554325775
$ws.Range("J4").Value = 0.02986826554325775
$ws.Range("M4").Value = 1.011383666666666
$ws.Range("N4").Value = 3.034151
$ws.Range("O4").Value = 0.9457805449469654
$ws.Range("P4").Value = 0.9457805449469653
$ws.Range("Q4").Value = 0.2592179708945555
$ws.Range("R4").Value = 2.332961738051
$ws.Range("S4").Value = 0.02824882446212298
$ws.Range("T4").Value = 0.02824882446212298

# Row 5
$ws.Range("G5").Value = 6.495645000000001
$ws.Range("I5").Value = 0.7569777503270297
$ws.Range("J5").Value = 0.7569777503270296
$ws.Range("K5").Value = 1.0
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04936366666666667
$ws.Range("N5").Value = 0.148091
$ws.Range("O5").Value = 0.04616170608573571
$ws.Range("P5").Value = 0.0461617060857357
$ws.Range("Q5").Value = 0.320648854565
$ws.Range("R5").Value = 2.885839691085001
$ws.Range("S5").Value = 0.03494338442403777
$ws.Range("T5").Value = 0.03494338442403776

# Row 6
$ws.Range("G6").Value = 6.495645000000001
$ws.Range("I6").Value = 0.7569777503270297
$ws.Range("J6").Value = 0.7569777503270296
$ws.Range("O6").Value = 0.008057748967298944
$ws.Range("P6").Value = 0.008057748967298944
$ws.Range("R6").Value = 0.5037372697500001
$ws.Range("S6").Value = 0.006099536685965902
$ws.Range("T6").Value = 0.006099536685965901

# Row 7
$ws.Range("G7").Value = 6.495645000000001
$ws.Range("I7").Value = 0.7569777503270297
$ws.Range("J7").Value = 0.7569777503270296
$ws.Range("M7").Value = 1.011383666666666
$ws.Range("N7").Value = 3.034151
$ws.Range("O7").Value = 0.9457805449469654
$ws.Range("P7").Value = 0.9457805449469653
$ws.Range("Q7").Value = 6.569589257464999
$ws.Range("R7").Value = 59.126303317185
$ws.Range("S7").Value = 0.7159348292170261
$ws.Range("T7").Value = 0.7159348292170259

# Row 8
$ws.Range("G8").Value = 1.804372666666667
$ws.Range("H8").Value = 5.413118000000001
$ws.Range("I8").Value = 0.2102747243676212
$ws.Range("J8").Value = 0.2102747243676212
$ws.Range("K8").Value = 1.0
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04936366666666667
$ws.Range("N8").Value = 0.148091
$ws.Range("O8").Value = 0.04616170608573571
$ws.Range("P8").Value = 0.0461617060857357
$ws.Range("Q8").Value = 0.08907045085977779
$ws.Range("R8").Value = 0.8016340577380001
$ws.Range("S8").Value = 0.009706640023517218
$ws.Range("T8").Value = 0.009706640023517215

# Row 9
$ws.Range("G9").Value = 1.804372666666667
$ws.Range("H9").Value = 5.413118000000001
$ws.Range("I9").Value = 0.2102747243676212
$ws.Range("J9").Value = 0.2102747243676212
$ws.Range("O9").Value = 0.008057748967298944
$ws.Range("P9").Value = 0.008057748967298944
$ws.Range("Q9").Value = 0.01554767781111111
$ws.Range("R9").Value = 0.1399291003
$ws.Range("S9").Value = 0.00169434094312227
$ws.Range("T9").Value = 0.00169434094312227

# Row 10
$ws.Range("G10").Value = 1.804372666666667
$ws.Range("H10").Value = 5.413118000000001
$ws.Range("I10").Value = 0.2102747243676212
$ws.Range("J10").Value = 0.2102747243676212
$ws.Range("M10").Value = 1.011383666666666
$ws.Range("N10").Value = 3.034151
$ws.Range("O10").Value = 0.9457805449469654
$ws.Range("P10").Value = 0.9457805449469653
$ws.Range("Q10").Value = 1.824913043646444
$ws.Range("R10").Value = 16.424217392818
$ws.Range("S10").Value = 0.1988737434009817
$ws.Range("T10").Value = 0.1988737434009817

# Row 11
$ws.Range("I11").Value = 0.002879259762091359
$ws.Range("J11").Value = 0.002879259762091358
$ws.Range("K11").Value = 1.0
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.04936366666666667
$ws.Range("N11").Value = 0.148091
$ws.Range("O11").Value = 0.04616170608573571
$ws.Range("P11").Value = 0.0461617060857357
$ws.Range("Q11").Value = 0.001219628112333334
$ws.Range("R11").Value = 0.010976653011
$ws.Range("S11").Value = 0.0001329115428821466
$ws.Range("T11").Value = 0.0001329115428821466

# Row 12
$ws.Range("I12").Value = 0.002879259762091359
$ws.Range("J12").Value = 0.002879259762091358
$ws.Range("O12").Value = 0.008057748967298944
$ws.Range("P12").Value = 0.008057748967298944
$ws.Range("S12").Value = 0.00002320035237457705
$ws.Range("T12").Value = 0.00002320035237457705

# Row 13
$ws.Range("I13").Value = 0.002879259762091359
$ws.Range("J13").Value = 0.002879259762091358
$ws.Range("M13").Value = 1.011383666666666
$ws.Range("N13").Value = 3.034151
$ws.Range("O13").Value = 0.9457805449469654
$ws.Range("P13").Value = 0.9457805449469653
$ws.Range("Q13").Value = 0.02498825625233333
$ws.Range("R13").Value = 0.224894306271
$ws.Range("S13").Value = 0.002723147866834635
$ws.Range("T13").Value = 0.002723147866834634
